$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 394, pushing the existing row 394 (and everything
# below it) down to rows 396+. This matches the diff, which shows the old
# rows 394-419 reappearing unchanged at 396-421, with two brand new rows of
# data landing at 394-395.
$ws.Rows.Item(394).Insert()
$ws.Rows.Item(394).Insert()

# New row 394
$ws.Cells.Item(394, 1).Value = 5
$ws.Cells.Item(394, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(394, 3).Value = 'Maule'
$ws.Cells.Item(394, 4).Value = 44826
$ws.Cells.Item(394, 5).Value = 7
$ws.Cells.Item(394, 6).Value = 100112023
$ws.Cells.Item(394, 7).Value = 'Brócoli'
$ws.Cells.Item(394, 8).Value = 'Sin especificar'
$ws.Cells.Item(394, 9).Value = 'Primera'
$ws.Cells.Item(394, 10).Value = 2000
$ws.Cells.Item(394, 11).Value = 1300
$ws.Cells.Item(394, 12).Value = 1300
$ws.Cells.Item(394, 13).Value = 1300
$ws.Cells.Item(394, 14).Value = '$/unidad'
$ws.Cells.Item(394, 15).Value = 'Región del Maule'
$ws.Cells.Item(394, 16).Value = 1300
$ws.Cells.Item(394, 17).Value = 1
$ws.Cells.Item(394, 18).Value = 'Hortaliza'

# New row 395
$ws.Cells.Item(395, 1).Value = 5
$ws.Cells.Item(395, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(395, 3).Value = 'Maule'
$ws.Cells.Item(395, 4).Value = 44826
$ws.Cells.Item(395, 5).Value = 7
$ws.Cells.Item(395, 6).Value = 100112023
$ws.Cells.Item(395, 7).Value = 'Brócoli'
$ws.Cells.Item(395, 8).Value = 'Sin especificar'
$ws.Cells.Item(395, 9).Value = 'Segunda'
$ws.Cells.Item(395, 10).Value = 2000
$ws.Cells.Item(395, 11).Value = 1000
$ws.Cells.Item(395, 12).Value = 1000
$ws.Cells.Item(395, 13).Value = 1000
$ws.Cells.Item(395, 14).Value = '$/unidad'
$ws.Cells.Item(395, 15).Value = 'Región del Maule'
$ws.Cells.Item(395, 16).Value = 1000
$ws.Cells.Item(395, 17).Value = 1
$ws.Cells.Item(395, 18).Value = 'Hortaliza'

Write-Host ("UsedRange after edit: " + $ws.UsedRange.Address())
